# Apply updated stock figures (quantities/values recalculated, plus a few
# swapped row-pairs and corrected totals) to the Companywise Stock Report.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B149").Value = 48654
$ws.Range("E149").Value = 38.26
$ws.Range("F149").Value = -1
$ws.Range("G149").Value = -32.02
$ws.Range("B150").Value = 63902
$ws.Range("E150").Value = 34.04
$ws.Range("F150").Value = 2
$ws.Range("G150").Value = 64.04000000000001
$ws.Range("B183").Value = 57552
$ws.Range("E183").Value = 136.86
$ws.Range("F183").Value = -5
$ws.Range("G183").Value = -603.45
$ws.Range("B184").Value = 64329
$ws.Range("E184").Value = 128.32
$ws.Range("F184").Value = 4
$ws.Range("G184").Value = 482.76
$ws.Range("F264").Value = 44
$ws.Range("G264").Value = 13013
$ws.Range("B310").Value = 132202.97
$ws.Range("B313").Value = 62997
$ws.Range("F313").Value = 0
$ws.Range("G313").Value = 0
$ws.Range("B314").Value = 57854
$ws.Range("F314").Value = 2
$ws.Range("G314").Value = 611.6799999999999
$ws.Range("B316").Value = 57077
$ws.Range("D316").Value = 93.08
$ws.Range("E316").Value = 111.2
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 93.08
$ws.Range("B317").Value = 63565
$ws.Range("D317").Value = 102.71
$ws.Range("E317").Value = 109.19
$ws.Range("F317").Value = 60
$ws.Range("G317").Value = 6162.6
$ws.Range("B350").Value = 57802
$ws.Range("E350").Value = 162.71
$ws.Range("F350").Value = -79
$ws.Range("G350").Value = -11334.92
$ws.Range("B351").Value = 63571
$ws.Range("E351").Value = 152.53
$ws.Range("F351").Value = 12
$ws.Range("G351").Value = 1721.76
$ws.Range("B375").Value = 63563
$ws.Range("E375").Value = 119.04
$ws.Range("F375").Value = 2
$ws.Range("G375").Value = 223.92
$ws.Range("B376").Value = 61605
$ws.Range("E376").Value = 133.78
$ws.Range("F376").Value = -13
$ws.Range("G376").Value = -1455.48
$ws.Range("B379").Value = 65514
$ws.Range("F379").Value = 0
$ws.Range("G379").Value = 0
$ws.Range("B380").Value = 63564
$ws.Range("F380").Value = 27
$ws.Range("G380").Value = 3483.27
$ws.Range("F387").Value = 4
$ws.Range("G387").Value = 234.8
$ws.Range("B389").Value = 57817
$ws.Range("F389").Value = 3
$ws.Range("G389").Value = 239.43
$ws.Range("B390").Value = 62865
$ws.Range("F390").Value = 4
$ws.Range("G390").Value = 319.24
$ws.Range("F420").Value = 822
$ws.Range("G420").Value = 140833.26
$ws.Range("B421").Value = 63008
$ws.Range("F421").Value = 421
$ws.Range("G421").Value = 63642.57
$ws.Range("B422").Value = 57857
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51
$ws.Range("B431").Value = 53082
$ws.Range("F431").Value = 1
$ws.Range("G431").Value = 59.47
$ws.Range("B432").Value = 63102
$ws.Range("F432").Value = 4
$ws.Range("G432").Value = 237.88
$ws.Range("B434").Value = 460093.97
$ws.Range("B536").Value = 47097
$ws.Range("D536").Value = 112.28
$ws.Range("E536").Value = 134.16
$ws.Range("F536").Value = 15
$ws.Range("G536").Value = 1684.2
$ws.Range("B537").Value = 58047
$ws.Range("D537").Value = 105.54
$ws.Range("E537").Value = 126.1
$ws.Range("F537").Value = 43
$ws.Range("G537").Value = 4538.22
$ws.Range("B583").Value = 53263
$ws.Range("E583").Value = 15.29
$ws.Range("F583").Value = -309
$ws.Range("G583").Value = -3958.29
$ws.Range("B584").Value = 65066
$ws.Range("E584").Value = 13.61
$ws.Range("F584").Value = 223
$ws.Range("G584").Value = 2856.63
$ws.Range("B586").Value = 64915
$ws.Range("E586").Value = 20.98
$ws.Range("F586").Value = 2
$ws.Range("G586").Value = 39.46
$ws.Range("B587").Value = 45695
$ws.Range("E587").Value = 23.58
$ws.Range("F587").Value = -36
$ws.Range("G587").Value = -710.28
$ws.Range("B590").Value = 64922
$ws.Range("E590").Value = 20.98
$ws.Range("F590").Value = 163
$ws.Range("G590").Value = 3215.99
$ws.Range("B591").Value = 45706
$ws.Range("E591").Value = 23.58
$ws.Range("F591").Value = -202
$ws.Range("G591").Value = -3985.46
$ws.Range("B593").Value = 45718
$ws.Range("E593").Value = 19.38
$ws.Range("F593").Value = -294
$ws.Range("G593").Value = -4768.68
$ws.Range("B594").Value = 64927
$ws.Range("E594").Value = 17.26
$ws.Range("F594").Value = 264
$ws.Range("G594").Value = 4282.08
$ws.Range("F628").Value = 104
$ws.Range("G628").Value = 1453.92
$ws.Range("F645").Value = 10
$ws.Range("G645").Value = 139.8
$ws.Range("B651").Value = 28632.84
$ws.Range("B687").Value = 53319
$ws.Range("E687").Value = 310.64
$ws.Range("F687").Value = -6
$ws.Range("G687").Value = -1643.52
$ws.Range("B688").Value = 64810
$ws.Range("E688").Value = 291.22
$ws.Range("F688").Value = 7
$ws.Range("G688").Value = 1917.44
$ws.Range("F815").Value = 0
$ws.Range("G815").Value = 0
$ws.Range("B831").Value = 60950.81
$ws.Range("B887").Value = 65362
$ws.Range("F887").Value = 77
$ws.Range("G887").Value = 3146.99
$ws.Range("B888").Value = 65079
$ws.Range("F888").Value = 21
$ws.Range("G888").Value = 858.27
$ws.Range("F932").Value = 3931
$ws.Range("G932").Value = 641185.41
$ws.Range("F934").Value = 663
$ws.Range("G934").Value = 187542.81
$ws.Range("F935").Value = 563
$ws.Range("G935").Value = 81437.95
$ws.Range("F937").Value = 17
$ws.Range("G937").Value = 1311.38
$ws.Range("F939").Value = 0
$ws.Range("G939").Value = 0
$ws.Range("B941").Value = 964550.1800000001
$ws.Range("B960").Value = 4153743.4
$ws.Range("B961").Value = 4153743.4
$ws.Range("C431").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("C432").Value = "HUL-Vim Bar Multipack Fw 4X200G"
